$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date value (45188) for all data rows
# (2-74). Update that date to 45189 (one day later) for each row.
for ($r = 2; $r -le 74; $r++) {
    $ws.Cells.Item($r, 3).Value = 45189
}
